# Update player-stats cells on Sheet1.
# All touched columns (E..L) are stored as text ("inlineStr") in the
# workbook, e.g. E3 holds the text "134", not the number 134. A leading
# apostrophe forces Excel to keep/store the new value as text too,
# instead of silently re-typing the cell as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "'224"
$ws.Range("F3").Value = "'3"
$ws.Range("G3").Value = "'2"

$ws.Range("J5").Value = "'2"

$ws.Range("E7").Value = "'1080"
$ws.Range("F7").Value = "'12"
$ws.Range("G7").Value = "'12"

$ws.Range("E8").Value = "'923"
$ws.Range("F8").Value = "'13"
$ws.Range("G8").Value = "'9"

$ws.Range("J9").Value = "'8"

$ws.Range("E10").Value = "'597"
$ws.Range("F10").Value = "'9"
$ws.Range("G10").Value = "'7"

$ws.Range("E11").Value = "'206"
$ws.Range("F11").Value = "'4"
$ws.Range("H11").Value = "'2"
$ws.Range("J11").Value = "'5"

$ws.Range("E12").Value = "'1149"
$ws.Range("F12").Value = "'13"
$ws.Range("G12").Value = "'13"
$ws.Range("I12").Value = "'1"
$ws.Range("L12").Value = "'3"

$ws.Range("E14").Value = "'104"
$ws.Range("F14").Value = "'3"
$ws.Range("G14").Value = "'1"

$ws.Range("J15").Value = "'9"

$ws.Range("E16").Value = "'691"
$ws.Range("F16").Value = "'13"
$ws.Range("G16").Value = "'9"
$ws.Range("I16").Value = "'8"

$ws.Range("E18").Value = "'787"
$ws.Range("F18").Value = "'12"
$ws.Range("G18").Value = "'10"

$ws.Range("J20").Value = "'10"

$ws.Range("E23").Value = "'255"
$ws.Range("F23").Value = "'8"
$ws.Range("G23").Value = "'3"
$ws.Range("I23").Value = "'3"

$ws.Range("E25").Value = "'1005"
$ws.Range("F25").Value = "'12"
$ws.Range("G25").Value = "'12"

$ws.Range("E26").Value = "'991"
$ws.Range("F26").Value = "'13"
$ws.Range("H26").Value = "'2"
$ws.Range("J26").Value = "'2"

$ws.Range("E27").Value = "'741"
$ws.Range("F27").Value = "'13"
$ws.Range("H27").Value = "'4"
$ws.Range("J27").Value = "'4"

$ws.Range("J28").Value = "'3"

$ws.Range("E29").Value = "'333"
$ws.Range("F29").Value = "'11"
$ws.Range("H29").Value = "'9"
$ws.Range("J29").Value = "'10"

$ws.Range("E31").Value = "'155"
$ws.Range("F31").Value = "'3"
$ws.Range("G31").Value = "'2"
$ws.Range("I31").Value = "'2"
